# Generate Report for Archive
# - Flip the localization status from "Ready for handoff" to "In Translation"
#   everywhere it appears (Overview!E2:F2, zh-cn!C2, de-de!C2 all share the
#   same shared-string entry).
# - Re-fit the Status/zh-cn/de-de columns that held that text so their
#   widths reflect the new (shorter) string.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $firstAddress = $null
    $found = $used.Find("Ready for handoff")
    while ($found -ne $null) {
        if ($firstAddress -eq $null) { $firstAddress = $found.Address }
        $found.Value = "In Translation"
        $found = $used.FindNext($found)
        if ($found -eq $null -or $found.Address -eq $firstAddress) { break }
    }
}

# Column width after the text got shorter ("Ready for handoff" -> "In
# Translation"): the Overview sheet keeps the status in columns E (zh-cn)
# and F (de-de); the per-locale sheets keep it in column C (Status).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E:F").ColumnWidth = 12.5

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C:C").ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C:C").ColumnWidth = 12.5
